$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

function Set-TextValue($cell, $value) {
    $cell.NumberFormat = "@"
    $cell.Value = $value
    $cell.Style = "Normal"
}

Set-TextValue $ws.Range("D2") "264.96"
Set-TextValue $ws.Range("G2") "6"
Set-TextValue $ws.Range("D3") "22.77"
Set-TextValue $ws.Range("G3") "6"
Set-TextValue $ws.Range("D4") "6.282"
Set-TextValue $ws.Range("G4") "6"
Set-TextValue $ws.Range("D5") "0.06152"
Set-TextValue $ws.Range("G5") "6"
Set-TextValue $ws.Range("D6") "3.592"
Set-TextValue $ws.Range("G6") "6"
Set-TextValue $ws.Range("D7") "6.697"
Set-TextValue $ws.Range("G7") "6"
Set-TextValue $ws.Range("D8") "1.345"
Set-TextValue $ws.Range("G8") "6"
Set-TextValue $ws.Range("D9") "0.8290"
Set-TextValue $ws.Range("G9") "6"
Set-TextValue $ws.Range("D10") "0.01356"
Set-TextValue $ws.Range("G10") "6"
Set-TextValue $ws.Range("G11") "6"
Set-TextValue $ws.Range("D12") "0.08242"
Set-TextValue $ws.Range("G12") "6"
Set-TextValue $ws.Range("D13") "0.03419"
Set-TextValue $ws.Range("G13") "6"
Set-TextValue $ws.Range("D14") "0.03135"
Set-TextValue $ws.Range("G14") "6"
Set-TextValue $ws.Range("D15") "0.09244"
Set-TextValue $ws.Range("G15") "6"
Set-TextValue $ws.Range("D16") "3.912"
Set-TextValue $ws.Range("G16") "6"
Set-TextValue $ws.Range("D17") "0.001725"
Set-TextValue $ws.Range("G17") "6"
Set-TextValue $ws.Range("D18") "0.04893"
Set-TextValue $ws.Range("G18") "6"
Set-TextValue $ws.Range("D19") "0.006307"
Set-TextValue $ws.Range("G19") "6"
Set-TextValue $ws.Range("D20") "0.005270"
Set-TextValue $ws.Range("G20") "6"
Set-TextValue $ws.Range("D21") "0.001089"
Set-TextValue $ws.Range("G21") "6"
Set-TextValue $ws.Range("G22") "6"
Set-TextValue $ws.Range("D23") "3.766"
Set-TextValue $ws.Range("G23") "6"
Set-TextValue $ws.Range("D24") "2.287"
Set-TextValue $ws.Range("G24") "6"
Set-TextValue $ws.Range("D25") "0.3380"
Set-TextValue $ws.Range("G25") "6"
Set-TextValue $ws.Range("G26") "6"
Set-TextValue $ws.Range("G27") "6"
Set-TextValue $ws.Range("G28") "6"
Set-TextValue $ws.Range("G29") "6"
Set-TextValue $ws.Range("G30") "6"
Set-TextValue $ws.Range("G31") "6"
Set-TextValue $ws.Range("G32") "6"
Set-TextValue $ws.Range("G33") "6"
Set-TextValue $ws.Range("G34") "6"
Set-TextValue $ws.Range("G35") "6"
Set-TextValue $ws.Range("G36") "6"
Set-TextValue $ws.Range("G37") "6"
Set-TextValue $ws.Range("G38") "6"
Set-TextValue $ws.Range("G39") "6"
Set-TextValue $ws.Range("D40") "0.04620"
Set-TextValue $ws.Range("G40") "6"
Set-TextValue $ws.Range("D41") "0.006955"
Set-TextValue $ws.Range("G41") "6"
Set-TextValue $ws.Range("D42") "0.1136"
Set-TextValue $ws.Range("G42") "6"
Set-TextValue $ws.Range("D43") "0.003130"
Set-TextValue $ws.Range("G43") "6"
Set-TextValue $ws.Range("D44") "0.01052"
Set-TextValue $ws.Range("G44") "6"
Set-TextValue $ws.Range("D45") "0.00006151"
Set-TextValue $ws.Range("G45") "6"
Set-TextValue $ws.Range("G46") "6"
Set-TextValue $ws.Range("G47") "6"
Set-TextValue $ws.Range("D48") "0.1964"
Set-TextValue $ws.Range("G48") "6"
Set-TextValue $ws.Range("G49") "6"
Set-TextValue $ws.Range("G50") "6"
Set-TextValue $ws.Range("G51") "6"
